# Update LeveProfits-style sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# with refreshed market-board price snapshots (scheduled data-refresh run).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 726.3684
$ws.Range("I15").Value = 726.3684
$ws.Range("K15").Value = 2179.1052
$ws.Range("M15").Value = -2010.1052
$ws.Range("H70").Value = 1346.1666
$ws.Range("J70").Value = 1443.0769
$ws.Range("L70").Value = 4329.2307
$ws.Range("N70").Value = -4869.2307
$ws.Range("H73").Value = 1346.1666
$ws.Range("J73").Value = 1443.0769
$ws.Range("L73").Value = 4329.2307
$ws.Range("N73").Value = -6201.2307
$ws.Range("H116").Value = 2090
$ws.Range("I116").Value = 1702.8572
$ws.Range("J116").Value = 4800
$ws.Range("K116").Value = 1702.8572
$ws.Range("L116").Value = 4800
$ws.Range("M116").Value = 1739.1428
$ws.Range("N116").Value = -11684
$ws.Range("H137").Value = 1270.4
$ws.Range("J137").Value = 1887.3
$ws.Range("L137").Value = 5661.9
$ws.Range("N137").Value = -10761.9
$ws.Range("H140").Value = 52925
$ws.Range("J140").Value = 52925
$ws.Range("L140").Value = 52925
$ws.Range("N140").Value = -63285

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7063.6914
$ws.Range("I32").Value = 8036.018
$ws.Range("J32").Value = 2950
$ws.Range("K32").Value = 8036.018
$ws.Range("L32").Value = 2950
$ws.Range("M32").Value = -7749.018
$ws.Range("N32").Value = -3524
$ws.Range("H74").Value = 5173.8335
$ws.Range("I74").Value = 3612.75
$ws.Range("K74").Value = 3612.75
$ws.Range("M74").Value = -2738.75
$ws.Range("H77").Value = 5173.8335
$ws.Range("I77").Value = 3612.75
$ws.Range("K77").Value = 18063.75
$ws.Range("M77").Value = -13695.75
$ws.Range("H92").Value = 51274.832
$ws.Range("I92").Value = 50000
$ws.Range("J92").Value = 51529.8
$ws.Range("K92").Value = 50000
$ws.Range("L92").Value = 51529.8
$ws.Range("M92").Value = -47504
$ws.Range("N92").Value = -56521.8
$ws.Range("H94").Value = 24776.666
$ws.Range("J94").Value = 24776.666
$ws.Range("L94").Value = 24776.666
$ws.Range("N94").Value = -26578.666
$ws.Range("H101").Value = 52602
$ws.Range("J101").Value = 52602
$ws.Range("L101").Value = 52602
$ws.Range("N101").Value = -59092
$ws.Range("H132").Value = 8623393
$ws.Range("I132").Value = 9617751
$ws.Range("J132").Value = 5621
$ws.Range("K132").Value = 28853253
$ws.Range("L132").Value = 16863
$ws.Range("M132").Value = -28850723
$ws.Range("N132").Value = -21923

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4340.3335
$ws.Range("I107").Value = 4340.3335
$ws.Range("K107").Value = 4340.3335
$ws.Range("M107").Value = -2420.3335
$ws.Range("H138").Value = 50800
$ws.Range("J138").Value = 50800
$ws.Range("L138").Value = 50800
$ws.Range("N138").Value = -61080

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5444.6206
$ws.Range("I31").Value = 7646.6313
$ws.Range("K31").Value = 7646.6313
$ws.Range("M31").Value = -7351.6313
$ws.Range("H34").Value = 5444.6206
$ws.Range("I34").Value = 7646.6313
$ws.Range("K34").Value = 7646.6313
$ws.Range("M34").Value = -7444.6313
$ws.Range("H58").Value = 2045
$ws.Range("I58").Value = 621.4286
$ws.Range("J58").Value = 5366.6665
$ws.Range("K58").Value = 621.4286
$ws.Range("L58").Value = 5366.6665
$ws.Range("M58").Value = -418.4286
$ws.Range("N58").Value = -5772.6665
$ws.Range("H94").Value = 5271.467
$ws.Range("I94").Value = 3453.6667
$ws.Range("J94").Value = 6483.3335
$ws.Range("K94").Value = 3453.6667
$ws.Range("L94").Value = 6483.3335
$ws.Range("M94").Value = -3002.6667
$ws.Range("N94").Value = -7385.3335
$ws.Range("H95").Value = 28333.334
$ws.Range("J95").Value = 28333.334
$ws.Range("L95").Value = 28333.334
$ws.Range("N95").Value = -33825.334
$ws.Range("H122").Value = 1528.4348
$ws.Range("I122").Value = 1556.4706
$ws.Range("J122").Value = 1449
$ws.Range("K122").Value = 4669.4118
$ws.Range("L122").Value = 4347
$ws.Range("M122").Value = -2219.4118
$ws.Range("N122").Value = -9247
$ws.Range("H132").Value = 17629.375
$ws.Range("I132").Value = 29006
$ws.Range("J132").Value = 6252.75
$ws.Range("K132").Value = 87018
$ws.Range("L132").Value = 18758.25
$ws.Range("M132").Value = -84488
$ws.Range("N132").Value = -23818.25
$ws.Range("H136").Value = 2045
$ws.Range("I136").Value = 621.4286
$ws.Range("J136").Value = 5366.6665
$ws.Range("K136").Value = 1864.2858
$ws.Range("L136").Value = 16099.9995
$ws.Range("M136").Value = 685.7142000000001
$ws.Range("N136").Value = -21199.9995

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 88.083336
$ws.Range("I8").Value = 88.083336
$ws.Range("K8").Value = 264.250008
$ws.Range("M8").Value = -125.250008
$ws.Range("H23").Value = 308
$ws.Range("J23").Value = 362.5
$ws.Range("L23").Value = 1087.5
$ws.Range("N23").Value = -1557.5
$ws.Range("H105").Value = 7343
$ws.Range("J105").Value = 7343
$ws.Range("L105").Value = 22029
$ws.Range("N105").Value = -27271
$ws.Range("H122").Value = 1468
$ws.Range("I122").Value = 1248.5834
$ws.Range("K122").Value = 11237.2506
$ws.Range("M122").Value = -8787.250599999999
$ws.Range("H132").Value = 1100.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1100.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9904.5
$ws.Range("N132").Value = -14964.5
$ws.Range("M132").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 9372.5
$ws.Range("J92").Value = 9372.5
$ws.Range("L92").Value = 9372.5
$ws.Range("N92").Value = -13116.5
$ws.Range("H122").Value = 1794.7646
$ws.Range("I122").Value = 1631
$ws.Range("J122").Value = 2327
$ws.Range("K122").Value = 4893
$ws.Range("L122").Value = 6981
$ws.Range("M122").Value = -2443
$ws.Range("N122").Value = -11881
$ws.Range("H126").Value = 4916.9
$ws.Range("I126").Value = 3077.2856
$ws.Range("J126").Value = 5907.4614
$ws.Range("K126").Value = 9231.856800000001
$ws.Range("L126").Value = 17722.3842
$ws.Range("M126").Value = -6761.856800000001
$ws.Range("N126").Value = -22662.3842
$ws.Range("H132").Value = 6421.2
$ws.Range("I132").Value = 5539.8
$ws.Range("K132").Value = 16619.4
$ws.Range("M132").Value = -14089.4
$ws.Range("H138").Value = 58313.715
$ws.Range("J138").Value = 58313.715
$ws.Range("L138").Value = 58313.715
$ws.Range("N138").Value = -68593.715

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 40083.75
$ws.Range("J94").Value = 40083.75
$ws.Range("L94").Value = 40083.75
$ws.Range("N94").Value = -41435.75
$ws.Range("H95").Value = 6000
$ws.Range("J95").Value = 6000
$ws.Range("L95").Value = 6000
$ws.Range("N95").Value = -11492
$ws.Range("H122").Value = 4830.3335
$ws.Range("I122").Value = 5784.5
$ws.Range("J122").Value = 4067
$ws.Range("K122").Value = 17353.5
$ws.Range("L122").Value = 12201
$ws.Range("M122").Value = -14903.5
$ws.Range("N122").Value = -17101
$ws.Range("H132").Value = 15160646
$ws.Range("I132").Value = 4994.2354
$ws.Range("J132").Value = 31263526
$ws.Range("K132").Value = 14982.7062
$ws.Range("L132").Value = 93790578
$ws.Range("M132").Value = -12452.7062
$ws.Range("N132").Value = -93795638
$ws.Range("H139").Value = 46409.4
$ws.Range("J139").Value = 47049.332
$ws.Range("L139").Value = 47049.332
$ws.Range("N139").Value = -57329.332

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 703
$ws.Range("I107").Value = 1716.2858
$ws.Range("J107").Value = 259.6875
$ws.Range("K107").Value = 5148.857400000001
$ws.Range("L107").Value = 779.0625
$ws.Range("M107").Value = -3228.857400000001
$ws.Range("N107").Value = -4619.0625
$ws.Range("H126").Value = 2524.6086
$ws.Range("I126").Value = 1760.2858
$ws.Range("J126").Value = 10550
$ws.Range("K126").Value = 5280.857400000001
$ws.Range("L126").Value = 31650
$ws.Range("M126").Value = -2810.857400000001
$ws.Range("N126").Value = -36590
$ws.Range("H132").Value = 4181.4546
$ws.Range("I132").Value = 2932
$ws.Range("J132").Value = 5222.6665
$ws.Range("K132").Value = 8796
$ws.Range("L132").Value = 15667.9995
$ws.Range("M132").Value = -6266
$ws.Range("N132").Value = -20727.9995
